$d = $word.ActiveDocument

# Find the last paragraph of the document ("Kommentare überarbeitet")
$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Last
$lastRange = $lastPara.Range

# Move to the very end of that paragraph's text (before the paragraph mark)
$insertPoint = $d.Range($lastRange.End - 1, $lastRange.End - 1)

$newTexts = @(
    "Buttons geändert (keine Variable mehr)",
    "Klassenvariablen teilweise zu lokalen Variablen geändert",
    "Fehlermeldungen bearbeitet"
)

foreach ($t in $newTexts) {
    $insertPoint.InsertParagraphAfter()
    $insertPoint = $d.Paragraphs.Last.Range
    $insertPoint.Text = $t
    $insertPoint.Collapse(0)
}
